# Applies the "合肥-漫展信息.xlsx" update:
#  - 展览(sheet1) & 全部类型(sheet4): bump "想去人数" (F) counters on a handful of
#    existing rows, and rename the row-20 event to include "第六届".
#  - 演出(sheet2): gains a brand-new event row (row 2) for the 包河留声机音乐节.
#  - 全部类型(sheet4) additionally gains the same new event as its new row 21
#    (it aggregates every category).
#  - 本地生活(sheet3) is untouched.

$wb = $excel.ActiveWorkbook

function Update-EventSheet($ws) {
    $ws.Range("F2").Value  = 8452
    $ws.Range("F3").Value  = 8104
    $ws.Range("F10").Value = 192
    $ws.Range("F11").Value = 243
    $ws.Range("F12").Value = 732
    $ws.Range("F14").Value = 3502
    $ws.Range("F19").Value = 144
    $ws.Range("F20").Value = 98
    $ws.Range("C20").Value = '合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~'
}

function Add-NewEvent($ws, [int]$row, [int]$idx) {
    # Column A carries the bold/centered/bordered "index" style (same as the
    # header + every other data row) - clone it from A1 before writing the value.
    $ws.Range("A1").Copy() | Out-Null
    $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("A" + $row).Value = $idx

    # "2024-08-03" parses as a date via COM automation, same as real Excel
    # would; force it back to literal text (matching the source workbook,
    # where every date column is stored as plain text) and strip the
    # leftover number-format style so the cell stays unstyled like its peers.
    $ws.Range("B" + $row).NumberFormat = "@"
    $ws.Range("B" + $row).Value = "2024-08-03"
    $ws.Range("B" + $row).ClearFormats()

    $ws.Range("C" + $row).Value = '合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会'
    $ws.Range("D" + $row).Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
    $ws.Range("E" + $row).Value = '2024.08.03 19:30-08.03 21:00'
    $ws.Range("F" + $row).Value = 0
    $ws.Range("G" + $row).Value = 80
    $ws.Range("H" + $row).Value = 'https://show.bilibili.com/platform/detail.html?id=83556'
    $ws.Range("I" + $row).Value = '//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg'
}

$wsExhibition = $wb.Worksheets.Item(1)  # 展览
$wsPerform    = $wb.Worksheets.Item(2)  # 演出
$wsAll        = $wb.Worksheets.Item(4)  # 全部类型

Update-EventSheet $wsExhibition
Update-EventSheet $wsAll

Add-NewEvent $wsPerform 2 1
Add-NewEvent $wsAll 21 20
